# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''311.63'
$ws.Range('E2').Value = '''0.73%'
$ws.Range('G2').Value = '''17'
$ws.Range('D3').Value = '''37.65'
$ws.Range('E3').Value = '''-0.14%'
$ws.Range('G3').Value = '''17'
$ws.Range('D4').Value = '''5.134'
$ws.Range('E4').Value = '''0.66%'
$ws.Range('G4').Value = '''17'
$ws.Range('D5').Value = '''0.07895'
$ws.Range('E5').Value = '''0.49%'
$ws.Range('G5').Value = '''17'
$ws.Range('D6').Value = '''4.411'
$ws.Range('E6').Value = '''0.88%'
$ws.Range('G6').Value = '''17'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D7').Value = '''8.267'
$ws.Range('E7').Value = '''-0.30%'
$ws.Range('G7').Value = '''17'
$ws.Range('B8').Value = 'FTXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D8').Value = '''1.896'
$ws.Range('E8').Value = '''-3.67%'
$ws.Range('G8').Value = '''17'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = '''0.9279'
$ws.Range('E9').Value = '''0.03%'
$ws.Range('G9').Value = '''17'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = '''0.1229'
$ws.Range('E10').Value = '''-9.19%'
$ws.Range('G10').Value = '''17'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '''0.1926'
$ws.Range('E11').Value = '''-3.97%'
$ws.Range('G11').Value = '''17'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '''0.09116'
$ws.Range('E12').Value = '''2.01%'
$ws.Range('G12').Value = '''17'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '''0.03308'
$ws.Range('E13').Value = '''-4.63%'
$ws.Range('G13').Value = '''17'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '''0.09632'
$ws.Range('E14').Value = '''-1.12%'
$ws.Range('G14').Value = '''17'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '''0.001380'
$ws.Range('E15').Value = '''-0.92%'
$ws.Range('G15').Value = '''17'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '''0.005692'
$ws.Range('E16').Value = '''-4.15%'
$ws.Range('G16').Value = '''17'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '''3.518'
$ws.Range('E17').Value = '''-1.84%'
$ws.Range('G17').Value = '''17'
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').Value = '''3.096'
$ws.Range('E18').Value = '''-1.19%'
$ws.Range('G18').Value = '''17'
$ws.Range('D19').Value = '''0.3399'
$ws.Range('E19').Value = '''-1.93%'
$ws.Range('G19').Value = '''17'
$ws.Range('D20').Value = '''5.301'
$ws.Range('E20').Value = '''5.81%'
$ws.Range('G20').Value = '''17'
$ws.Range('D21').Value = '''0.1272'
$ws.Range('E21').Value = '''-1.72%'
$ws.Range('G21').Value = '''17'
$ws.Range('D22').Value = '''0.2591'
$ws.Range('E22').Value = '''3.07%'
$ws.Range('G22').Value = '''17'
$ws.Range('G23').Value = '''17'
$ws.Range('E24').Value = '''1.16%'
$ws.Range('G24').Value = '''17'
$ws.Range('E25').Value = '''1.61%'
$ws.Range('G25').Value = '''17'
$ws.Range('E26').Value = '''-5.21%'
$ws.Range('G26').Value = '''17'
$ws.Range('D27').Value = '''0.0001221'
$ws.Range('E27').Value = '''-9.72%'
$ws.Range('G27').Value = '''17'
$ws.Range('G28').Value = '''17'
$ws.Range('G29').Value = '''17'
$ws.Range('G30').Value = '''17'
$ws.Range('G31').Value = '''17'
$ws.Range('G32').Value = '''17'
$ws.Range('G33').Value = '''17'
$ws.Range('G34').Value = '''17'
$ws.Range('G35').Value = '''17'
$ws.Range('G36').Value = '''17'
$ws.Range('G37').Value = '''17'
$ws.Range('G38').Value = '''17'
$ws.Range('D39').Value = '''0.02107'
$ws.Range('E39').Value = '''-8.40%'
$ws.Range('G39').Value = '''17'
$ws.Range('D40').Value = '''0.05166'
$ws.Range('E40').Value = '''2.20%'
$ws.Range('G40').Value = '''17'
$ws.Range('D41').Value = '''0.007582'
$ws.Range('E41').Value = '''1.49%'
$ws.Range('G41').Value = '''17'
$ws.Range('D42').Value = '''0.009143'
$ws.Range('E42').Value = '''-7.27%'
$ws.Range('G42').Value = '''17'
$ws.Range('E43').Value = '''0.34%'
$ws.Range('G43').Value = '''17'
$ws.Range('D44').Value = '''0.002061'
$ws.Range('E44').Value = '''3.96%'
$ws.Range('G44').Value = '''17'
$ws.Range('D45').Value = '''0.008621'
$ws.Range('E45').Value = '''-1.72%'
$ws.Range('G45').Value = '''17'
$ws.Range('D46').Value = '''0.00006707'
$ws.Range('E46').Value = '''-1.41%'
$ws.Range('G46').Value = '''17'
$ws.Range('E47').Value = '''-0.10%'
$ws.Range('G47').Value = '''17'
$ws.Range('D48').Value = '''0.001201'
$ws.Range('E48').Value = '''-7.78%'
$ws.Range('G48').Value = '''17'
$ws.Range('D49').Value = '''0.002806'
$ws.Range('E49').Value = '''-6.61%'
$ws.Range('G49').Value = '''17'
$ws.Range('D50').Value = '''0.00002101'
$ws.Range('E50').Value = '''-0.10%'
$ws.Range('G50').Value = '''17'
$ws.Range('D51').Value = '''0.0002001'
$ws.Range('E51').Value = '''-0.10%'
$ws.Range('G51').Value = '''17'
